{"js": "// \"minor update on STAR reports\"\n//\n// 1) \". This is a \" -> \". This is \" in the Scaled Score paragraph.\n// 2) The (cursor-position) \"_GoBack\" bookmark moves from the end of the\n//    \"...try doing some of these activities at home: \" paragraph to sit\n//    right before the \"ss_cg\" merge-field placeholder in the Scaled Score\n//    paragraph (i.e. immediately after the edited \". This is \" run).\n\nconst body = context.document.body;\n\n// --- 1. Fix the wording: \". This is a \" -> \". This is \" -------------------\nconst wordingMatches = body.search(\". This is a \", { matchCase: true });\nwordingMatches.load(\"text\");\nawait context.sync();\n\nif (wordingMatches.items.length > 0) {\n  wordingMatches.items[0].insertText(\". This is \", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- 2. Relocate the \"_GoBack\" bookmark ------------------------------------\n// Remove it from its old spot (if present there / anywhere).\nconst oldBookmark = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\noldBookmark.load(\"isNullObject\");\nawait context.sync();\n\nif (!oldBookmark.isNullObject) {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// Re-insert it, collapsed, immediately before the \"ss_cg\" placeholder text\n// that now directly follows the \". This is \" run.\nconst ssCgMatches = body.search(\"ss_cg\", { matchCase: true });\nssCgMatches.load(\"text\");\nawait context.sync();\n\nif (ssCgMatches.items.length > 0) {\n  const insertionPoint = ssCgMatches.items[0].getRange(Word.RangeLocation.start);\n  insertionPoint.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# \"minor update on STAR reports\"\n#\n# 1) \". This is a \" -> \". This is \" in the Scaled Score paragraph.\n# 2) The (cursor-position) \"_GoBack\" bookmark moves from the end of the\n#    \"...try doing some of these activities at home: \" paragraph to sit\n#    right before the \"ss_cg\" merge-field placeholder in the Scaled Score\n#    paragraph (i.e. immediately after the edited \". This is \" run).\n\n$d = $word.ActiveDocument\n\n# --- 1. Fix the wording: \". This is a \" -> \". This is \" -------------------\n$wordingRange = $d.Content\n$wordingRange.Find.Execute(\". This is a \", $false, $false, $false, $false, $false, $true, 1, $false, \". This is \", 2)\n\n# --- 2. Relocate the \"_GoBack\" bookmark ------------------------------------\n# Remove it from its old spot, if it is there.\ntry {\n    $d.Bookmarks(\"_GoBack\").Delete()\n} catch {\n}\n\n# Re-insert it, collapsed, immediately before the \"ss_cg\" placeholder text\n# that now directly follows the \". This is \" run.\n$bmRange = $d.Content\nif ($bmRange.Find.Execute(\"ss_cg\")) {\n    $insertPoint = $d.Range($bmRange.Start, $bmRange.Start)\n    $d.Bookmarks.Add(\"_GoBack\", $insertPoint)\n}\n"}
